$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.662.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.17%  "

$ws.Range("D3").Value = "'3.791.02"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.21%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'594.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("D6").Value = "'166.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.86%  "

$ws.Range("D7").Value = "'3.791.12"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.26%  "

$ws.Range("E8").Value = "  +0.12%  "

$ws.Range("E9").Value = "  -0.39%  "

$ws.Range("E10").Value = "  -0.62%  "

$ws.Range("E11").Value = "  -1.97%  "

$ws.Range("D12").Value = "'0.450"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").Value = "'0.0000255"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.80%  "

$ws.Range("D14").Value = "'36.22"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.47%  "

$ws.Range("D15").Value = "'4.425.15"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.15%  "

$ws.Range("D16").Value = "'3.787.14"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.53%  "

$ws.Range("D17").Value = "'18.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.97%  "

$ws.Range("D18").Value = "'67.654.54"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.13%  "

$ws.Range("E19").Value = "  +0.39%  "

$ws.Range("E20").Value = "  -0.99%  "

$ws.Range("E21").Value = "  -5.03%  "

$ws.Range("D22").Value = "'457.71"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.13%  "

$ws.Range("D23").Value = "'0.698"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.78%  "

$ws.Range("D24").Value = "'0.0000151"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.06%  "

$ws.Range("D25").Value = "'83.64"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.73%  "

$ws.Range("D26").Value = "'11.86"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.13%  "

$ws.Range("E27").Value = "  -3.11%  "

$ws.Range("D28").Value = "'10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.61%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("D30").Value = "'2.79"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.25%  "

$ws.Range("D31").Value = "'7.29"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.86%  "

$ws.Range("D32").Value = "'29.89"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.25%  "

$ws.Range("D33").Value = "'2.19"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.01%  "

$ws.Range("D34").Value = "'9.19"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.24%  "

$ws.Range("E35").Value = "  -0.35%  "

$ws.Range("D36").Value = "'3.743.98"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.17%  "

$ws.Range("E37").Value = "  -1.55%  "

$ws.Range("D38").Value = "'3.33"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.50%  "

$ws.Range("E39").Value = "  -0.25%  "

$ws.Range("D40").Value = "'0.995"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("E41").Value = "  -1.28%  "

$ws.Range("E42").Value = "  +0.08%  "

$ws.Range("E43").Value = "  +0.01%  "

$ws.Range("D44").Value = "'44.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.39%  "

$ws.Range("E45").Value = "  -2.69%  "

$ws.Range("D46").Value = "'47.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.73%  "

$ws.Range("D47").Value = "'8.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.91%  "

$ws.Range("D48").Value = "'147.83"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.78%  "

$ws.Range("D49").Value = "'391.39"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.96%  "

$ws.Range("E50").Value = "  -5.78%  "

$ws.Range("D51").Value = "'2.753.81"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.41%  "
